$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds dates stored as serial 45233 for rows 2-89.
# Update them all to 45243 (10 day shift), preserving existing formatting.
$ws.Range("C2:C89").Value = 45243
